# Nerul.xlsx update
# Sheet3 ("Petrol/Site log" style sheet, dimension A1:E50) gets:
#  - several date corrections in column B (rows 14-20)
#  - newly filled-in data rows 21-26 (previously blank placeholder rows)
#  - a refreshed selection / scroll position in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Column B date corrections (rows 14-20): values shift forward by 30 days ---
$ws.Range("B14").Value = 45202
$ws.Range("B15").Value = 45203
$ws.Range("B16").Value = 45203
$ws.Range("B17").Value = 45203
$ws.Range("B18").Value = 45203
$ws.Range("B19").Value = 45203
$ws.Range("B20").Value = 45204

# --- Row 21 ---
$ws.Range("B20").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B21").Value = 45205

$ws.Range("C20").Copy()
$ws.Range("C21").PasteSpecial(-4122)
$ws.Range("C21").Value = "Deepak"

$ws.Range("D20").Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("D21").Value = "Flexible"

$ws.Range("E21").Value = 170

# --- Row 22 ---
$ws.Range("B20").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B22").Value = 45205

$ws.Range("C20").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = "Deepak"

$ws.Range("D20").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = "Petrol (18094)"

$ws.Range("E22").Value = 100

# --- Row 23 ---
$ws.Range("B20").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B23").Value = 45206

$ws.Range("C20").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = "Sandesh"

$ws.Range("D20").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D23").Value = "Petrol (No Bill)"

$ws.Range("E23").Value = 200

# --- Row 24 ---
$ws.Range("B20").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("B24").Value = 45208

$ws.Range("C20").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("C24").Value = "Deepak"

$ws.Range("D20").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("D24").Value = "Petrol (No Bill)"

$ws.Range("E24").Value = 200

# --- Row 25 (C/E keep the plain border style; D gets a brand-new centered, borderless style) ---
$ws.Range("B20").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("B25").Value = 45208

$ws.Range("C25").Value = "Sir"

$d25 = $ws.Range("D25")
$d25.Borders.LineStyle = -4142
$d25.VerticalAlignment = -4107
$d25.HorizontalAlignment = -4108
$d25.WrapText = $true
$d25.Font.Size = 10
$d25.Value = "Fabricator Nerul"

$ws.Range("E25").Value = 800

# --- Row 26 ---
$ws.Range("B20").Copy()
$ws.Range("B26").PasteSpecial(-4122)
$ws.Range("B26").Value = 45208

$ws.Range("C26").Value = "Deepak"
$ws.Range("D26").Value = "Rikshaw"
$ws.Range("E26").Value = 700

# --- Refresh the sheet view's selection to match the author's last editing position ---
$ws.Activate()
$ws.Range("E27").Select()
